# Tsalka Municipality.xlsx — collapse the 1989/2002/2014 area table down to
# just the 2014 figure, and drop the "(according to the population census
# data)" note row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "(according to the population census data)" note in A2 — it is
# dropped entirely in the new layout.
$ws.Range("A2").Clear()

# Drop the blank spacer row (old row 3); everything below shifts up one row.
$ws.Range("A3").EntireRow.Delete()

# Remove the 1989 and 2002 columns (B:C); the 2014 column (D) slides left
# into B.
$ws.Range("B4:C4").EntireColumn.Delete()

# Re-apply the row heights used by the simplified layout.
$ws.Range("A1:B6").RowHeight = 20.1

$wb.Save()
